$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.198.31'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = '2.418.27'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("E4").Value = '  +0.07%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '563.06'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.00%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '142.98'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +1.52%  '
$ws.Range("D9").Value = '2.414.10'
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("E10").Value = '  +1.30%  '
$ws.Range("E11").Value = '  -2.22%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '5.33'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.35%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.353'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.18%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '25.68'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D16").Value = '2.853.72'
$ws.Range("E16").Value = '  +1.22%  '
$ws.Range("D17").Value = '62.034.92'
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("D18").Value = '2.413.44'
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("E19").Value = '  +1.87%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '4.18'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.37%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '323.46'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.22%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.85'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +2.88%  '
$ws.Range("E23").Value = '  -0.18%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '65.78'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +2.02%  '
$ws.Range("E25").Value = '  -1.79%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '8.99'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.13%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '578.92'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +4.42%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0949'
$ws.Range("E28").Value = '  +3.37%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '2.531.17'
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.13%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '8.24'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.49%  '
$ws.Range("E32").Value = '  +0.92%  '
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("E34").Value = '  +1.07%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  +0.01%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.56'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -3.48%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '4.73'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.25%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '152.10'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +4.36%  '
$ws.Range("E41").Value = '  +0.43%  '
$ws.Range("E42").Value = '  -7.50%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.995'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("E44").Value = '  +1.42%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '148.61'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("E46").Value = '  +0.55%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.0535'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.16%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '20.06'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("E49").Value = '  +1.50%  '
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("E51").Value = '  +1.66%  '
